# Client_EPIS_Daily_Progress.xlsx - Team_Management sheet: remove one row
# (the data previously on row 11 is gone; rows 12-21 shift up to 11-20),
# shrink the filtered range/used dimension from G21 to G20 accordingly,
# and leave the selection where Excel lands after the delete (D12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team_Management")

# Delete row 11 entirely - remaining rows below shift up and carry their
# original formatting/styles with them, matching Excel's normal row-delete
# behavior.
[void]$ws.Rows.Item(11).Delete()

# The sheet's AutoFilter range doesn't auto-shrink when rows are deleted,
# so toggle it off and reapply over the new (smaller) data range.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
[void]$ws.Range("A1:G20").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase defined name in sync with
# the resized AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Team_Management!_FilterDatabase") {
        $n.RefersTo = "=Team_Management!`$A`$1:`$G`$20"
    }
}

# Leave the selection where it ended up after the edit.
[void]$ws.Range("D12").Select()
